$wb = $excel.ActiveWorkbook

# --- Pelamar: mark Haji Ahda's interview status as finished ---
$pelamar = $wb.Worksheets.Item("Pelamar")
$pelamar.Range("E2").Value = "Selesai"

# --- Wawancara: add new interview schedule row for Haji Ahda (Proggrammer) ---
# Tanggal/jam kept as plain text (matches existing rows), so force the
# "@" text format before assigning to stop Excel auto-converting them to
# a date serial / number.
$wawancara = $wb.Worksheets.Item("Wawancara")
$wawancara.Range("A3:F3").NumberFormat = "@"
$wawancara.Cells.Item(3, 1).Value = "L003"
$wawancara.Cells.Item(3, 2).Value = "P001"
$wawancara.Cells.Item(3, 3).Value = "Haji Ahda"
$wawancara.Cells.Item(3, 4).Value = "Proggrammer"
$wawancara.Cells.Item(3, 5).Value = "2024-01-20"
$wawancara.Cells.Item(3, 6).Value = "16.00"

# --- Add new "Seleksi" sheet right after "Wawancara" ---
$seleksi = $wb.Worksheets.Add($null, $wawancara)
$seleksi.Name = "Seleksi"

$seleksi.Cells.Item(1, 1).Value = "Kode Pelamar"
$seleksi.Cells.Item(1, 2).Value = "Nama"
$seleksi.Cells.Item(1, 3).Value = "Posisi"
$seleksi.Cells.Item(1, 4).Value = "Status Seleksi"

$seleksi.Cells.Item(2, 1).Value = "P001"
$seleksi.Cells.Item(2, 2).Value = "Haji Ahda"
$seleksi.Cells.Item(2, 3).Value = "Proggrammer"
$seleksi.Cells.Item(2, 4).Value = "Diterima"

# Restore the original active sheet (Pelamar) since adding a sheet
# activates it by default.
$pelamar.Activate()
